$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 724 (everything from the old row 724 onward
# shifts down by two, matching the diff: old rows 724-765 become 726-767).
$ws.Rows.Item(724).Insert()
$ws.Rows.Item(724).Insert()

# New row 724: 2026/01/25, 日, 23, 148
$ws.Range("A724").NumberFormat = "@"
$ws.Range("A724").Value = "2026/01/25"
$ws.Range("A724").ClearFormats()
$ws.Range("B724").Value = "日"
$ws.Range("C724").Value = 23
$ws.Range("D724").Value = 148

# New row 725: 2026/01/26, 月, 2, 158
$ws.Range("A725").NumberFormat = "@"
$ws.Range("A725").Value = "2026/01/26"
$ws.Range("A725").ClearFormats()
$ws.Range("B725").Value = "月"
$ws.Range("C725").Value = 2
$ws.Range("D725").Value = 158
